# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "22.405.23"
$ws.Range("E2").Value = "  -3.99%  "
# Row 3
$ws.Range("D3").Value = "1.575.03"
$ws.Range("E3").Value = "  -3.26%  "
# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.21%  "
# Row 5
$ws.Range("E5").Value = "  -0.15%  "
# Row 6
$ws.Range("D6").Value = "'289.91"
$ws.Range("E6").Value = "  -2.79%  "
# Row 7
$ws.Range("D7").Value = "'0.3680"
$ws.Range("E7").Value = "  -2.37%  "
# Row 8
$ws.Range("D8").Value = "'49.39"
$ws.Range("E8").Value = "  -1.40%  "
# Row 9
$ws.Range("D9").Value = "'0.3387"
$ws.Range("E9").Value = "  -3.45%  "
# Row 10
$ws.Range("E10").Value = "  -2.55%  "
# Row 11
$ws.Range("D11").Value = "'0.07626"
# Row 12
$ws.Range("E12").Value = "  -0.18%  "
# Row 13
$ws.Range("D13").Value = "'21.37"
$ws.Range("E13").Value = "  -2.06%  "
# Row 14
$ws.Range("E14").Value = "  -3.34%  "
# Row 15
$ws.Range("D15").Value = "'6.936"
$ws.Range("E15").Value = "  -3.79%  "
# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.579.20"
$ws.Range("E16").Value = "  -3.11%  "
# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.00001139"
$ws.Range("E17").Value = "  -4.28%  "
# Row 18
$ws.Range("E18").Value = "  -6.15%  "
# Row 19
$ws.Range("D19").Value = "'0.06741"
$ws.Range("E19").Value = "  -2.85%  "
# Row 20
$ws.Range("E20").Value = "  -0.22%  "
# Row 21
$ws.Range("D21").Value = "'6.255"
$ws.Range("E21").Value = "  -6.04%  "
# Row 22
$ws.Range("D22").Value = "'16.60"
$ws.Range("E22").Value = "  -3.71%  "
# Row 23
$ws.Range("D23").Value = "'0.5328"
$ws.Range("E23").Value = "  -6.60%  "
# Row 24
$ws.Range("D24").Value = "'11.99"
$ws.Range("E24").Value = "  -1.93%  "
# Row 25
$ws.Range("D25").Value = "22.422.95"
$ws.Range("E25").Value = "  -4.03%  "
# Row 26
$ws.Range("E26").Value = "  -2.62%  "
# Row 27
$ws.Range("D27").Value = "'2.992"
$ws.Range("E27").Value = "  +2.57%  "
# Row 28
$ws.Range("D28").Value = "'20.03"
$ws.Range("E28").Value = "  -3.22%  "
# Row 29
$ws.Range("D29").Value = "'145.74"
$ws.Range("E29").Value = "  -3.92%  "
# Row 30
$ws.Range("D30").Value = "'4.986"
$ws.Range("E30").Value = "  -3.38%  "
# Row 32
$ws.Range("D32").Value = "1.752.32"
$ws.Range("E32").Value = "  -3.30%  "
# Row 33
$ws.Range("D33").Value = "'1.048"
$ws.Range("E33").Value = "  +8.63%  "
# Row 34
$ws.Range("D34").Value = "'6.303"
$ws.Range("E34").Value = "  -6.83%  "
# Row 35
$ws.Range("D35").Value = "'2.000"
$ws.Range("E35").Value = "  -5.48%  "
# Row 36
$ws.Range("D36").Value = "'10.36"
$ws.Range("E36").Value = "  -7.48%  "
# Row 37
$ws.Range("D37").Value = "'0.08456"
# Row 38
$ws.Range("D38").Value = "'0.02544"
$ws.Range("E38").Value = "  -5.18%  "
# Row 39
$ws.Range("D39").Value = "'0.2331"
$ws.Range("E39").Value = "  -3.63%  "
# Row 40
$ws.Range("D40").Value = "'0.06590"
$ws.Range("E40").Value = "  -2.52%  "
# Row 41
$ws.Range("D41").Value = "'5.558"
$ws.Range("E41").Value = "  -4.47%  "
# Row 42
$ws.Range("D42").Value = "'11.82"
$ws.Range("E42").Value = "  -7.20%  "
# Row 43
$ws.Range("E43").Value = "  -3.22%  "
# Row 44
$ws.Range("D44").Value = "'0.6390"
$ws.Range("E44").Value = "  -5.90%  "
# Row 45
$ws.Range("D45").Value = "'14.47"
$ws.Range("E45").Value = "  -5.82%  "
# Row 46
$ws.Range("D46").Value = "'0.9998"
$ws.Range("E46").Value = "  -0.22%  "
# Row 47
$ws.Range("D47").Value = "'0.6010"
$ws.Range("E47").Value = "  -4.45%  "
# Row 48
$ws.Range("E48").Value = "  -3.55%  "
# Row 49
$ws.Range("D49").Value = "'2.132"
$ws.Range("E49").Value = "  -4.05%  "
# Row 50
$ws.Range("D50").Value = "'1.260"
$ws.Range("E50").Value = "  +5.46%  "
# Row 51
$ws.Range("D51").Value = "'123.52"
$ws.Range("E51").Value = "  -2.26%  "
